# "added new params for evpcd_i"
#
# Summary of changes applied:
#  - evp-cd_f (sheet2): mergeCells reordered to canonical order; pane
#    scroll/selection moved; 4 existing data rows re-labelled (C column
#    text changes only).
#  - evp-cd_i (sheet3): mergeCells reordered to canonical order; pane
#    scroll/selection moved; 8 new parameter rows populated (11,12,23,24,
#    35,36,47,48) that were previously blank placeholders.
#  - evp-wd_f (sheet4) / evp-wd_i (sheet5): mergeCells reordered only.

$wb = $excel.ActiveWorkbook

function Reorder-Merges {
    param($ws, [string[]]$allRanges, [string[]]$targetOrder)

    foreach ($r in $allRanges) {
        if ($ws.Range($r).MergeCells) {
            $ws.Range($r).UnMerge()
        }
    }
    foreach ($r in $targetOrder) {
        $ws.Range($r).Merge()
    }
}

$mergeRanges = @("A2:A13","B2:B13","A14:A25","B14:B25","A26:A37","B26:B37","A38:A49","B38:B49","A50:A61","B50:B61","A62:A73","B62:B73")
$orderFrontGroup  = @("A2:A13","B2:B13","A14:A25","B14:B25","A26:A37","B26:B37","A38:A49","B38:B49","A50:A61","B50:B61","A62:A73","B62:B73")
$orderBackGroup   = @("A38:A49","B38:B49","A50:A61","B50:B61","A62:A73","B62:B73","A2:A13","B2:B13","A14:A25","B14:B25","A26:A37","B26:B37")

# ---------------------------------------------------------------
# evp-cd_f  (sheet2)
# ---------------------------------------------------------------
$wsCdF = $wb.Worksheets.Item("evp-cd_f")

Reorder-Merges -ws $wsCdF -allRanges $mergeRanges -targetOrder $orderFrontGroup

# Re-labelled data points (only the "C" text label changes; the measured
# values in D:K are untouched)
$wsCdF.Cells.Item(24, 3).Value = "9; 9"
$wsCdF.Cells.Item(35, 3).Value = "8; 8"
$wsCdF.Cells.Item(36, 3).Value = "8; 8"
$wsCdF.Cells.Item(48, 3).Value = "8; 8"

$wsCdF.Activate()
$wsCdF.Range("C48").Select()

# ---------------------------------------------------------------
# evp-cd_i  (sheet3) - the sheet the commit message calls out
# ---------------------------------------------------------------
$wsCdI = $wb.Worksheets.Item("evp-cd_i")

Reorder-Merges -ws $wsCdI -allRanges $mergeRanges -targetOrder $orderBackGroup

# Row 11
$wsCdI.Cells.Item(11, 3).Value = "9; 7"
$wsCdI.Cells.Item(11, 4).Value = 19.125
$wsCdI.Cells.Item(11, 5).Value = 43.641
$wsCdI.Cells.Item(11, 6).Value = 5.6148
$wsCdI.Cells.Item(11, 7).Value = 4.1688
$wsCdI.Cells.Item(11, 8).Value = 1616
$wsCdI.Cells.Item(11, 9).Value = 1876.8
$wsCdI.Cells.Item(11, 10).Value = 5.5594
$wsCdI.Cells.Item(11, 11).Value = 6.8653

# Row 12
$wsCdI.Cells.Item(12, 3).Value = "9; 8"
$wsCdI.Cells.Item(12, 4).Value = 22.012
$wsCdI.Cells.Item(12, 5).Value = 38.873
$wsCdI.Cells.Item(12, 6).Value = 2.2033
$wsCdI.Cells.Item(12, 7).Value = 4.245
$wsCdI.Cells.Item(12, 8).Value = 1841.2
$wsCdI.Cells.Item(12, 9).Value = 2377
$wsCdI.Cells.Item(12, 10).Value = 5.2898
$wsCdI.Cells.Item(12, 11).Value = 9.7712

# Row 23
$wsCdI.Cells.Item(23, 3).Value = "9; 7"
$wsCdI.Cells.Item(23, 4).Value = 5.0569
$wsCdI.Cells.Item(23, 5).Value = 40.476
$wsCdI.Cells.Item(23, 6).Value = 10.017
$wsCdI.Cells.Item(23, 7).Value = 4.1585
$wsCdI.Cells.Item(23, 8).Value = 1730.1
$wsCdI.Cells.Item(23, 9).Value = 1998.1
$wsCdI.Cells.Item(23, 10).Value = 5.5564
$wsCdI.Cells.Item(23, 11).Value = 10.337

# Row 24
$wsCdI.Cells.Item(24, 3).Value = "9; 9"
$wsCdI.Cells.Item(24, 4).Value = 0.36499
$wsCdI.Cells.Item(24, 5).Value = 45.374
$wsCdI.Cells.Item(24, 6).Value = 15.687
$wsCdI.Cells.Item(24, 7).Value = 3.6105
$wsCdI.Cells.Item(24, 8).Value = 2878.3
$wsCdI.Cells.Item(24, 9).Value = 3325.8
$wsCdI.Cells.Item(24, 10).Value = 4.847
$wsCdI.Cells.Item(24, 11).Value = 12.328

# Row 35
$wsCdI.Cells.Item(35, 3).Value = "9; 10"
$wsCdI.Cells.Item(35, 4).Value = 8.0268
$wsCdI.Cells.Item(35, 5).Value = 28.899
$wsCdI.Cells.Item(35, 6).Value = 2.5811
$wsCdI.Cells.Item(35, 7).Value = 4.018
$wsCdI.Cells.Item(35, 8).Value = 1213.4
$wsCdI.Cells.Item(35, 9).Value = 1698.5
$wsCdI.Cells.Item(35, 10).Value = 4.5711
$wsCdI.Cells.Item(35, 11).Value = 7.3221

# Row 36
$wsCdI.Cells.Item(36, 3).Value = "9; 8"
$wsCdI.Cells.Item(36, 4).Value = 8.0846
$wsCdI.Cells.Item(36, 5).Value = 21.623
$wsCdI.Cells.Item(36, 6).Value = 4.9985
$wsCdI.Cells.Item(36, 7).Value = 3.8734
$wsCdI.Cells.Item(36, 8).Value = 1330.8
$wsCdI.Cells.Item(36, 9).Value = 2586.3
$wsCdI.Cells.Item(36, 10).Value = 4.0925
$wsCdI.Cells.Item(36, 11).Value = 6.8218

# Row 47
$wsCdI.Cells.Item(47, 3).Value = "9; 9"
$wsCdI.Cells.Item(47, 4).Value = 11.849
$wsCdI.Cells.Item(47, 5).Value = 272.72
$wsCdI.Cells.Item(47, 6).Value = 0.21347
$wsCdI.Cells.Item(47, 7).Value = 3.8447
$wsCdI.Cells.Item(47, 8).Value = 1238.3
$wsCdI.Cells.Item(47, 9).Value = 1731
$wsCdI.Cells.Item(47, 10).Value = 4.5635
$wsCdI.Cells.Item(47, 11).Value = 7.4758

# Row 48
$wsCdI.Cells.Item(48, 3).Value = "8; 9"
$wsCdI.Cells.Item(48, 4).Value = 7.9997
$wsCdI.Cells.Item(48, 5).Value = 25.008
$wsCdI.Cells.Item(48, 6).Value = 2.2597
$wsCdI.Cells.Item(48, 7).Value = 4.127
$wsCdI.Cells.Item(48, 8).Value = 1102.2
$wsCdI.Cells.Item(48, 9).Value = 1841.7
$wsCdI.Cells.Item(48, 10).Value = 4.5535
$wsCdI.Cells.Item(48, 11).Value = 10.206

$wsCdI.Activate()
$wsCdI.Range("I23").Select()

# ---------------------------------------------------------------
# evp-wd_f  (sheet4) - mergeCells order only
# ---------------------------------------------------------------
$wsWdF = $wb.Worksheets.Item("evp-wd_f")
Reorder-Merges -ws $wsWdF -allRanges $mergeRanges -targetOrder $orderFrontGroup

# ---------------------------------------------------------------
# evp-wd_i  (sheet5) - mergeCells order only
# ---------------------------------------------------------------
$wsWdI = $wb.Worksheets.Item("evp-wd_i")
Reorder-Merges -ws $wsWdI -allRanges $mergeRanges -targetOrder $orderBackGroup

# Restore the sheet that was active/selected in the source workbook.
$wsCdI.Activate()
